$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price text like "25.841.47" that Excel would otherwise
# auto-parse as a number; force Text format first so values round-trip as
# strings exactly like the source workbook.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '25.849.47'
$ws.Range('E2').Value = '  +0.34%  '
$ws.Range('D3').Value = '1.637.59'
$ws.Range('E3').Value = '  +0.64%  '
$ws.Range('D4').Value = '1.005'
$ws.Range('E4').Value = '  +0.32%  '
$ws.Range('D5').Value = '215.29'
$ws.Range('E5').Value = '  -0.07%  '
$ws.Range('D6').Value = '0.5093'
$ws.Range('E6').Value = '  -0.30%  '
$ws.Range('E7').Value = '  +0.24%  '
$ws.Range('D8').Value = '0.2585'
$ws.Range('E8').Value = '  +0.93%  '
$ws.Range('D9').Value = '0.06439'
$ws.Range('E9').Value = '  +1.87%  '
$ws.Range('D10').Value = '20.42'
$ws.Range('E10').Value = '  +5.01%  '
$ws.Range('D11').Value = '0.07797'
$ws.Range('E11').Value = '  +0.23%  '
$ws.Range('B12').Value = 'Polkadot'
$ws.Range('C12').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D12').Value = '4.260'
$ws.Range('E12').Value = '  +0.65%  '
$ws.Range('B13').Value = 'WrappedEther'
$ws.Range('C13').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D13').Value = '1.638.69'
$ws.Range('E13').Value = '  +0.22%  '
$ws.Range('D14').Value = '1.862.28'
$ws.Range('E14').Value = '  +0.58%  '
$ws.Range('D15').Value = '0.5605'
$ws.Range('E15').Value = '  +1.59%  '
$ws.Range('D16').Value = '0.0₅7680'
$ws.Range('E16').Value = '  +1.74%  '
$ws.Range('D17').Value = '63.31'
$ws.Range('E17').Value = '  -0.53%  '
$ws.Range('D18').Value = '25.863.96'
$ws.Range('E18').Value = '  +0.31%  '
$ws.Range('D19').Value = '1.005'
$ws.Range('E19').Value = '  +0.28%  '
$ws.Range('D20').Value = '4.385'
$ws.Range('E20').Value = '  -0.46%  '
$ws.Range('D21').Value = '192.65'
$ws.Range('E21').Value = '  -0.75%  '
$ws.Range('D22').Value = '9.985'
$ws.Range('E22').Value = '  +1.50%  '
$ws.Range('D23').Value = '6.152'
$ws.Range('E23').Value = '  +2.51%  '
$ws.Range('E24').Value = '  +0.23%  '
$ws.Range('D25').Value = '1.757'
$ws.Range('E25').Value = '  -6.82%  '
$ws.Range('D26').Value = '139.29'
$ws.Range('E26').Value = '  -1.89%  '
$ws.Range('E27').Value = '  -1.92%  '
$ws.Range('D28').Value = '6.839'
$ws.Range('E28').Value = '  +1.44%  '
$ws.Range('E29').Value = '  -0.08%  '
$ws.Range('D30').Value = '1.241'
$ws.Range('E30').Value = '  +0.21%  '
$ws.Range('D31').Value = '0.04974'
$ws.Range('E31').Value = '  +1.91%  '
$ws.Range('D32').Value = '3.310'
$ws.Range('E32').Value = '  +2.56%  '
$ws.Range('D33').Value = '3.263'
$ws.Range('E33').Value = '  +2.92%  '
$ws.Range('D34').Value = '1.574'
$ws.Range('E34').Value = '  +2.35%  '
$ws.Range('D35').Value = '2.391'
$ws.Range('E35').Value = '  +0.60%  '
$ws.Range('D36').Value = '0.9041'
$ws.Range('E36').Value = '  +1.24%  '
$ws.Range('D37').Value = '0.5579'
$ws.Range('D38').Value = '2.573'
$ws.Range('E38').Value = '  +1.69%  '
$ws.Range('D39').Value = '1.134.77'
$ws.Range('E39').Value = '  +2.09%  '
$ws.Range('E40').Value = '  +1.41%  '
$ws.Range('D41').Value = '0.9969'
$ws.Range('E41').Value = '  -0.32%  '
$ws.Range('B42').Value = 'FraxShare'
$ws.Range('C42').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D42').Value = '5.462'
$ws.Range('E42').Value = '  -1.65%  '
$ws.Range('B43').Value = 'TrustWalletToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D43').Value = '0.8017'
$ws.Range('E43').Value = '  +0.63%  '
$ws.Range('B44').Value = 'Quant'
$ws.Range('C44').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D44').Value = '99.13'
$ws.Range('E44').Value = '  +2.03%  '
$ws.Range('D45').Value = '0.0₈114'
$ws.Range('E45').Value = '  +1.85%  '
$ws.Range('D46').Value = '55.78'
$ws.Range('E46').Value = '  +2.19%  '
$ws.Range('D47').Value = '0.4270'
$ws.Range('E47').Value = '  -3.59%  '
$ws.Range('D48').Value = '7.822'
$ws.Range('E48').Value = '  +3.59%  '
$ws.Range('D49').Value = '0.05032'
$ws.Range('E49').Value = '  -1.99%  '
$ws.Range('D50').Value = '0.9993'
$ws.Range('E50').Value = '  -0.27%  '
$ws.Range('D51').Value = '1.004'
$ws.Range('E51').Value = '  +0.43%  '
